$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-10: Sending cluster x Target cluster combinations (ECs, FAPs, sCs)^2
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$rows = @(
    @("ECs",  "Gnai2", "P2ry12", "ECs",  3, 1, 203.7816646666667, 611.344994,        0.6667327591988204,  0.6667327591988205,  1, 0.3333333333333333, 7.960095666666667, 23.880287, 0.8748886527059857, 0.8748886527059857, 1622.121545859253, 14599.09391273328, 0.5833169254104004,  0.5833169254104005),
    @("ECs",  "Gnai2", "P2ry12", "FAPs", 3, 1, 203.7816646666667, 611.344994,        0.6667327591988204,  0.6667327591988205,  1, 0.3333333333333333, 0.019266,          0.057798,  0.002117512840155588, 0.002117512840155588, 3.926057551468,    35.334517963212,   0.001411815178555866, 0.001411815178555866),
    @("ECs",  "Gnai2", "P2ry12", "sCs",  3, 1, 203.7816646666667, 611.344994,        0.6667327591988204,  0.6667327591988205,  3, 1,                  1.119048333333333, 3.357145,  0.1229938344538588,   0.1229938344538588,   228.0415322091256, 2052.37378988213,  0.08200401860986421, 0.08200401860986423),
    @("FAPs", "Gnai2", "P2ry12", "ECs",  3, 1, 63.14058933333333, 189.421768,        0.2065833519051582,  0.2065833519051582,  1, 0.3333333333333333, 7.960095666666667, 23.880287, 0.8748886527059857, 0.8748886527059857, 502.6051315430462, 4523.446183887416, 0.1807374304197903,  0.1807374304197904),
    @("FAPs", "Gnai2", "P2ry12", "FAPs", 3, 1, 63.14058933333333, 189.421768,        0.2065833519051582,  0.2065833519051582,  1, 0.3333333333333333, 0.019266,          0.057798,  0.002117512840155588, 0.002117512840155588, 1.216466594096,    10.948199346864,   0.0004374429002215528,0.0004374429002215529),
    @("FAPs", "Gnai2", "P2ry12", "sCs",  3, 1, 63.14058933333333, 189.421768,        0.2065833519051582,  0.2065833519051582,  3, 1,                  1.119048333333333, 3.357145,  0.1229938344538588,   0.1229938344538588,   70.65737125915112, 635.9163413323599, 0.02540847858514628, 0.02540847858514628),
    @("sCs",  "Gnai2", "P2ry12", "ECs",  3, 1, 38.719942,          116.159826,       0.1266838888960214,  0.1266838888960214,  1, 0.3333333333333333, 7.960095666666667, 23.880287, 0.8748886527059857, 0.8748886527059857, 308.2144425277846, 2773.929982750062, 0.1108342968757949,  0.110834296875795),
    @("sCs",  "Gnai2", "P2ry12", "FAPs", 3, 1, 38.719942,          116.159826,       0.1266838888960214,  0.1266838888960214,  1, 0.3333333333333333, 0.019266,          0.057798,  0.002117512840155588, 0.002117512840155588, 0.745978402572,    6.713805623148,    0.0002682547613781692,0.0002682547613781693),
    @("sCs",  "Gnai2", "P2ry12", "sCs",  3, 1, 38.719942,          116.159826,       0.1266838888960214,  0.1266838888960214,  3, 1,                  1.119048333333333, 3.357145,  0.1229938344538588,   0.1229938344538588,   43.32948656186333, 389.96537905677,   0.0155813372588483,  0.0155813372588483)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
